# "Enhance name parsing logic for multi-person names and suffixes"
#
# Updates the name-parsing test-data sheet so the expected/parsed values
# reflect the refined regex logic described in the commit:
#   - suffix extraction now keeps the shared last name together with the
#     suffix (row 4: "Jr" -> "Smith Jr", and the blank last-name cell now
#     captures "Smith")
#   - multi-person name handling now pulls the correct first name for the
#     *other* person sharing a last name (rows 5, 7, 23, 24)
#   - the "Last, First" / multi-word-surname case now reorders/reformats
#     the expected output (row 10)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "Mr. John H. Smith Jr. & Mrs. Jane Smith"
#   LastName (C4) was left blank by the old logic; the new logic recovers
#   the shared surname, and Full Lastname (D4) keeps the suffix attached.
$ws.Range("C4").Value = "Smith"
$ws.Range("D4").Value = "Smith Jr"

# Row 5: "Lydia & Gary Orange" -- FirstName should be the second person's
# first name, not the first.
$ws.Range("B5").Value = "Gary"

# Row 7: "Ben J. and Maggie E. Goddy" -- same fix, second person's first name.
$ws.Range("B7").Value = "Maggie"

# Row 10: "García-López, Maria" -- Expected Output reformatted to
# "<Last> <First>" with the accented surname preserved.
$ws.Range("E10").Value = "García-López Maria"

# Row 23: "Michael & Susan Thompson" -- FirstName corrected to "Susan".
$ws.Range("B23").Value = "Susan"

# Row 24: "David & Mary Wilson" -- FirstName corrected to "Mary".
$ws.Range("B24").Value = "Mary"
